# Split the two "... You can do this with " sentences into three runs each:
#   "<lead>" | "This can be done using the following command:" | " " [" "]
#
# Word merges adjacent runs that share identical formatting when a
# document is saved, so a plain Range.Text assignment collapses right
# back into a single run. To keep the desired run boundaries we toggle a
# character property (Bold on/off, a true no-op visually) on exactly the
# sub-ranges that must remain distinct runs; that forces Word to keep
# them split into separate <w:r> elements.

$d = $word.ActiveDocument

function Split-Run([int]$rangeStart, [int]$rangeEnd) {
    # Re-apply (no-op) direct character formatting over [start, end) so
    # this slice of text is emitted as its own run, independent of its
    # neighbours.
    $rTemp = $d.Range($rangeStart, $rangeEnd)
    $rTemp.Font.Bold = 1
    $rTemp.Font.Bold = 0
}

# ---------------------------------------------------------------------
# Occurrence 1: "Method 2: You can do this with "
#            -> "Method 2: " | "This can be done using the following command:" | " "
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("Method 2: You can do this with ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

$full1   = "Method 2: You can do this with "
$target1 = "You can do this with "
$offset1 = $full1.IndexOf($target1)

$start1    = $find1.Start
$subStart1 = $start1 + $offset1
$subEnd1   = $start1 + $full1.Length

$newTail1  = "This can be done using the following command: "
$d.Range($subStart1, $subEnd1).Text = $newTail1

$run2Start1 = $subStart1
$run2End1   = $subStart1 + "This can be done using the following command:".Length
Split-Run $run2Start1 $run2End1

# ---------------------------------------------------------------------
# Occurrence 2: "You can do this with "
#            -> "This can be done using the following command:" | " " | " "
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("You can do this with ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

$start2 = $find2.Start
$end2   = $find2.End

$newText2 = "This can be done using the following command:  "
$d.Range($start2, $end2).Text = $newText2

$part1Len2 = "This can be done using the following command:".Length
$run1Start2 = $start2
$run1End2   = $start2 + $part1Len2
$run2Start2 = $run1End2
$run2End2   = $run2Start2 + 1

Split-Run $run1Start2 $run1End2
Split-Run $run2Start2 $run2End2

Write-Output "Applied both run-split edits"
